# Refresh portfolio valuations (simulated market-data pull) and push the
# results into Portfolio, Daily_Summary and Performance_History.

$wb = $excel.ActiveWorkbook
$portfolio = $wb.Worksheets.Item("Portfolio")

# Row 2 (TATAMOTORS) - the price lookup failed this run, so the derived
# columns go blank and the daily change resets to 0.
$portfolio.Range("G2").Formula = '=""'
$portfolio.Range("H2").Formula = '=""'
$portfolio.Range("I2").Formula = '=""'
$portfolio.Range("J2").Formula = '=""'
$portfolio.Range("K2").Value = 0
$portfolio.Range("L2").Value = "2025-11-18 10:20:13"

# Rows 3-15: refreshed Current_Price (G) plus the dependent columns
# (Current_Value, Unrealized_PL, PL_Percentage, Daily_PL_Percentage) and the
# refresh timestamp.
$rows = @(
    @{ Row = 3;  G = 172.45;  H = 10002.1;            I = 589.2799999999988;   J = 6.26039805286831;   K = -0.4157764046890332 }
    @{ Row = 4;  G = 1486.4;  H = 1486.4;              I = -3.079999999999927;  J = -0.2067835754760002; K = -1.406208543380195 }
    @{ Row = 5;  G = 915.5;   H = 43944;               I = 1272;                J = 2.98087739032621;    K = 0.01638717430490821 }
    @{ Row = 6;  G = 510.7;   H = 7660.5;              I = 3977.25;             J = 107.9820810425575;   K = -1.939324116743465 }
    @{ Row = 7;  G = 381.5;   H = 12589.5;             I = -188.4299999999985;  J = -1.474651997624016;  K = -0.3135615364515257 }
    @{ Row = 8;  G = 386.25;  H = 15450;               I = 4997.6;              J = 47.81294248210938;   K = -1.654996817313813 }
    @{ Row = 9;  G = 28.38;   H = 737.88;              I = -210.08;             J = -22.16127262753703;  K = -0.2460456942003525 }
    @{ Row = 10; G = 246.95;  H = 1234.75;             I = 472.35;              J = 61.95566631689402;   K = -0.443458980044355 }
    @{ Row = 11; G = 826.7;   H = 4133.5;              I = -416.9499999999998;  J = -9.162830049775293;  K = -0.9287554676733176 }
    @{ Row = 12; G = 328.45;  H = 4926.75;             I = 59.54999999999927;   J = 1.223496055226809;   K = -0.5299818291944276 }
    @{ Row = 13; G = 703.8;   H = 3519;                I = -197.9000000000001;  J = -5.324329414296862;  K = -1.297244232522264 }
    @{ Row = 14; G = 136.65;  H = 10932;               I = -2220.799999999999;  J = -16.88461772398272;  K = -0.6543075245365362 }
    @{ Row = 15; G = 1595.2;  H = 4785.6;              I = 1688.820000000001;   J = 54.5347102474183;    K = -0.6972111553784888 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $portfolio.Range("G$row").Value = $r.G
    $portfolio.Range("H$row").Value = $r.H
    $portfolio.Range("I$row").Value = $r.I
    $portfolio.Range("J$row").Value = $r.J
    $portfolio.Range("K$row").Value = $r.K
    $portfolio.Range("L$row").Value = "2025-11-18 10:20:13"
}

# Daily_Summary and Performance_History both roll up the same totals for the
# new as-of date. The leading apostrophe keeps "2025-11-18" a literal text
# value instead of Excel auto-converting it to a date serial; re-applying the
# Normal style afterwards drops the quote-prefix formatting it would
# otherwise pick up.
$summarySheets = @("Daily_Summary", "Performance_History")
foreach ($name in $summarySheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("A2").Value = "'2025-11-18"
    $ws.Range("A2").Style = "Normal"
    $ws.Range("B2").Value = 116363.41
    $ws.Range("C2").Value = 121401.98
    $ws.Range("D2").Value = 5038.570000000022
    $ws.Range("E2").Value = 4.330029517010564
}
